$d = $word.ActiveDocument

# 1) "Chuyen nganh" paragraph: the phrase ", mã " (one run) + "chương trình"
#    (the following run) together read ", mã chương trình" before the colon.
#    Trim the first run down to just "," and empty out the second run so the
#    ": Kỹ thuật..." run is left immediately after the comma.
$d.Content.Find.Execute(", mã chương trình:", $true, $false, $false, $false, $false,
                         $true, 1, $false, ",:", 2) | Out-Null

# 2) "Chuong trinh dao tao" paragraph: insert ", mã chuyên ngành" before the
#    closing period that follows "mô tả về chương trình".
$d.Content.Find.Execute("mô tả về chương trình.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mô tả về chương trình, mã chuyên ngành.", 2) | Out-Null

# 3) Remove the _GoBack bookmark left over at the end of the document.
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}
